$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 821.0476
$ws.Range("I6").Value = 138.22223
$ws.Range("K6").Value = 414.66669
$ws.Range("M6").Value = -302.66669
$ws.Range("H112").Value = 1294.527
$ws.Range("I112").Value = 694
$ws.Range("J112").Value = 1338.0435
$ws.Range("K112").Value = 2082
$ws.Range("L112").Value = 4014.1305
$ws.Range("M112").Value = -974
$ws.Range("N112").Value = -6230.1305
$ws.Range("H129").Value = 1069.2361
$ws.Range("I129").Value = 495.5
$ws.Range("J129").Value = 1121.3939
$ws.Range("K129").Value = 1486.5
$ws.Range("L129").Value = 3364.1817
$ws.Range("M129").Value = 3513.5
$ws.Range("N129").Value = -13364.1817
$ws.Range("H132").Value = 2056.647
$ws.Range("I132").Value = 2056.647
$ws.Range("K132").Value = 6169.941
$ws.Range("M132").Value = -3639.941
$ws.Range("H138").Value = 5640.625
$ws.Range("I138").Value = 1144.0526
$ws.Range("J138").Value = 8586.655000000001
$ws.Range("K138").Value = 3432.1578
$ws.Range("L138").Value = 25759.965
$ws.Range("M138").Value = 1707.8422
$ws.Range("N138").Value = -36039.965

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2141051.5
$ws.Range("I122").Value = 2335328.8
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 7005986.399999999
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -7003536.399999999
$ws.Range("N122").Value = -16900
$ws.Range("H132").Value = 3140.5642
$ws.Range("I132").Value = 1470.7391
$ws.Range("J132").Value = 5540.9375
$ws.Range("K132").Value = 4412.2173
$ws.Range("L132").Value = 16622.8125
$ws.Range("M132").Value = -1882.2173
$ws.Range("N132").Value = -21682.8125

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3992.6123
$ws.Range("I134").Value = 4627.0884
$ws.Range("J134").Value = 2554.4666
$ws.Range("K134").Value = 13881.2652
$ws.Range("L134").Value = 7663.399800000001
$ws.Range("M134").Value = -11346.2652
$ws.Range("N134").Value = -12733.3998

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4591.625
$ws.Range("I31").Value = 1113.0625
$ws.Range("K31").Value = 1113.0625
$ws.Range("M31").Value = -818.0625
$ws.Range("H34").Value = 4591.625
$ws.Range("I34").Value = 1113.0625
$ws.Range("K34").Value = 1113.0625
$ws.Range("M34").Value = -911.0625

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 74.333336
$ws.Range("I14").Value = 74.333336
$ws.Range("K14").Value = 223.000008
$ws.Range("M14").Value = -50.00000800000001
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H68").Value = 1767.5
$ws.Range("I68").Value = 440
$ws.Range("K68").Value = 1320
$ws.Range("M68").Value = -509
$ws.Range("H71").Value = 1767.5
$ws.Range("I71").Value = 440
$ws.Range("K71").Value = 3960
$ws.Range("M71").Value = 96
$ws.Range("H86").Value = 920.55554
$ws.Range("I86").Value = 894.2857
$ws.Range("J86").Value = 1012.5
$ws.Range("K86").Value = 2682.8571
$ws.Range("L86").Value = 3037.5
$ws.Range("M86").Value = -1496.8571
$ws.Range("N86").Value = -5409.5
$ws.Range("H89").Value = 920.55554
$ws.Range("I89").Value = 894.2857
$ws.Range("J89").Value = 1012.5
$ws.Range("K89").Value = 8048.571300000001
$ws.Range("L89").Value = 9112.5
$ws.Range("M89").Value = -2120.571300000001
$ws.Range("N89").Value = -20968.5
$ws.Range("H92").Value = 592
$ws.Range("J92").Value = 900
$ws.Range("L92").Value = 2700
$ws.Range("N92").Value = -5196
$ws.Range("H113").Value = 556149.7
$ws.Range("I113").Value = 550
$ws.Range("J113").Value = 1111749.4
$ws.Range("K113").Value = 1650
$ws.Range("L113").Value = 3335248.2
$ws.Range("M113").Value = 520
$ws.Range("N113").Value = -3339588.2
$ws.Range("H114").Value = 26028.5
$ws.Range("I114").Value = 676
$ws.Range("J114").Value = 51381
$ws.Range("K114").Value = 2028
$ws.Range("L114").Value = 154143
$ws.Range("M114").Value = 1226
$ws.Range("N114").Value = -160651
$ws.Range("H129").Value = 2153.889
$ws.Range("I129").Value = 1962
$ws.Range("J129").Value = 2227.6924
$ws.Range("K129").Value = 5886
$ws.Range("L129").Value = 6683.0772
$ws.Range("M129").Value = -886
$ws.Range("N129").Value = -16683.0772
$ws.Range("H131").Value = 1667639
$ws.Range("I131").Value = 7143650
$ws.Range("J131").Value = 1027.0435
$ws.Range("K131").Value = 21430950
$ws.Range("L131").Value = 3081.1305
$ws.Range("M131").Value = -21425910
$ws.Range("N131").Value = -13161.1305

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4469.926
$ws.Range("I132").Value = 4630.909
$ws.Range("J132").Value = 4359.25
$ws.Range("K132").Value = 13892.727
$ws.Range("L132").Value = 13077.75
$ws.Range("M132").Value = -11362.727
$ws.Range("N132").Value = -18137.75
$ws.Range("H137").Value = 38120
$ws.Range("J137").Value = 38120
$ws.Range("L137").Value = 38120
$ws.Range("N137").Value = -48320

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 64225.188
$ws.Range("I7").Value = 78806.766
$ws.Range("J7").Value = 1038.3334
$ws.Range("K7").Value = 78806.766
$ws.Range("L7").Value = 1038.3334
$ws.Range("M7").Value = -78694.766
$ws.Range("N7").Value = -1262.3334
$ws.Range("H16").Value = 3500
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 3500
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 3500
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -3840
$ws.Range("H61").Value = 2931.2856
$ws.Range("I61").Value = 2603.8
$ws.Range("K61").Value = 2603.8
$ws.Range("M61").Value = -2401.8
$ws.Range("H113").Value = 2931.2856
$ws.Range("I113").Value = 2603.8
$ws.Range("K113").Value = 2603.8
$ws.Range("M113").Value = -433.8000000000002
$ws.Range("H126").Value = 64225.188
$ws.Range("I126").Value = 78806.766
$ws.Range("J126").Value = 1038.3334
$ws.Range("K126").Value = 236420.298
$ws.Range("L126").Value = 3115.0002
$ws.Range("M126").Value = -233950.298
$ws.Range("N126").Value = -8055.0002
$ws.Range("H136").Value = 5636.85
$ws.Range("I136").Value = 5639.32
$ws.Range("J136").Value = 5632.7334
$ws.Range("K136").Value = 16917.96
$ws.Range("L136").Value = 16898.2002
$ws.Range("M136").Value = -14367.96
$ws.Range("N136").Value = -21998.2002

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 28000
$ws.Range("J86").Value = 28000
$ws.Range("L86").Value = 28000
$ws.Range("N86").Value = -30246
$ws.Range("H89").Value = 28000
$ws.Range("J89").Value = 28000
$ws.Range("L89").Value = 140000
$ws.Range("N89").Value = -151232
